$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 18

# Read all the data rows into an array of hashtables
$data = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @{
        A = $ws.Cells.Item($r, 1).Value2
        B = $ws.Cells.Item($r, 2).Value2
        C = $ws.Cells.Item($r, 3).Value2
        D = $ws.Cells.Item($r, 4).Value2
    }
    $data += $row
}

# Sort ascending by column A (time)
$sorted = $data | Sort-Object { $_.A }

# Write the sorted data back
$r = $firstRow
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 1).Value2 = $row.A
    $ws.Cells.Item($r, 2).Value2 = $row.B
    $ws.Cells.Item($r, 3).Value2 = $row.C
    $ws.Cells.Item($r, 4).Value2 = $row.D
    $r++
}
